# Weekly update: insert two new price rows (newest date) at the top of the
# data block for "Agrícola del Norte S.A. de Arica - Acelga", pushing the
# existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 31 (existing rows 31:43 shift to 33:45).
$ws.Rows("31:32").Insert()

# New row 31: Acelga, "Primera" quality, newest report date.
$ws.Range("A31").Value = 1
$ws.Range("B31").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C31").Value = "Arica y Parinacota"
$ws.Range("D31").Value = 44523
$ws.Range("E31").Value = 15
$ws.Range("F31").Value = 100112009
$ws.Range("G31").Value = "Acelga"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 1400
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = 1450
$ws.Range("N31").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 483
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = "Hortaliza"

# New row 32: Acelga, "Segunda" quality, same newest report date.
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C32").Value = "Arica y Parinacota"
$ws.Range("D32").Value = 44523
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 100112009
$ws.Range("G32").Value = "Acelga"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 1100
$ws.Range("M32").Value = 1050
$ws.Range("N32").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 350
$ws.Range("Q32").Value = 3
$ws.Range("R32").Value = "Hortaliza"
